$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.346005439758301
$ws.Range("B1").Value = 4.979316711425781
$ws.Range("C1").Value = 6.32066011428833
$ws.Range("D1").Value = 10.4017333984375
$ws.Range("E1").Value = 3.59139609336853
